$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like text columns (Y, AA) to be stored as plain text instead
# of being auto-converted to Excel date serials, then restore the default
# "Normal" style so no stray number-format styling is left on the cells.
$dateCols = @(25, 27)  # Y, AA
foreach ($r in 3..5) {
    foreach ($c in $dateCols) {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
    }
}

# Columns that are present but hold an explicit empty string (Antal,
# Bestämningsår, Projektnamn). Assigning "" through .Value clears/omits the
# cell entirely in this host, so use a formula that evaluates to an empty
# string to materialize an actual (empty) text cell there instead.
$emptyTextCols = @(9, 46, 51)  # I, AT, AY

# Row 3
$ws.Cells.Item(3, 1).Value = 112551629
$ws.Cells.Item(3, 2).Value = 77650
$ws.Cells.Item(3, 3).Value = "Ovaliderad"
$ws.Cells.Item(3, 4).Value = "NT"
$ws.Cells.Item(3, 5).Value = 6425
$ws.Cells.Item(3, 6).Value = "Garnlav"
$ws.Cells.Item(3, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(3, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(3, 16).Value = "Syd Åskaksdammen, Vrm"
$ws.Cells.Item(3, 17).Value = 409494
$ws.Cells.Item(3, 18).Value = 6717338
$ws.Cells.Item(3, 19).Value = 10
$ws.Cells.Item(3, 20).Value = "Värmland"
$ws.Cells.Item(3, 21).Value = "Torsby"
$ws.Cells.Item(3, 22).Value = "Värmland"
$ws.Cells.Item(3, 23).Value = "Norra Ny"
$ws.Cells.Item(3, 25).Value = "2022-08-18"
$ws.Cells.Item(3, 27).Value = "2022-08-18"
$ws.Cells.Item(3, 30).Value = $false
$ws.Cells.Item(3, 31).Value = $false
$ws.Cells.Item(3, 33).Value = $false
$ws.Cells.Item(3, 49).Value = "Anders Boström"
$ws.Cells.Item(3, 50).Value = "Anders Boström"
foreach ($c in $emptyTextCols) {
    $ws.Cells.Item(3, $c).Formula = '=""'
}

# Row 4
$ws.Cells.Item(4, 1).Value = 112551631
$ws.Cells.Item(4, 2).Value = 77650
$ws.Cells.Item(4, 3).Value = "Ovaliderad"
$ws.Cells.Item(4, 4).Value = "NT"
$ws.Cells.Item(4, 5).Value = 6425
$ws.Cells.Item(4, 6).Value = "Garnlav"
$ws.Cells.Item(4, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(4, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(4, 16).Value = "Sydväst Åskaksdammen, Vrm"
$ws.Cells.Item(4, 17).Value = 409335
$ws.Cells.Item(4, 18).Value = 6717479
$ws.Cells.Item(4, 19).Value = 10
$ws.Cells.Item(4, 20).Value = "Värmland"
$ws.Cells.Item(4, 21).Value = "Torsby"
$ws.Cells.Item(4, 22).Value = "Värmland"
$ws.Cells.Item(4, 23).Value = "Norra Ny"
$ws.Cells.Item(4, 25).Value = "2022-08-18"
$ws.Cells.Item(4, 27).Value = "2022-08-18"
$ws.Cells.Item(4, 30).Value = $false
$ws.Cells.Item(4, 31).Value = $false
$ws.Cells.Item(4, 33).Value = $false
$ws.Cells.Item(4, 49).Value = "Anders Boström"
$ws.Cells.Item(4, 50).Value = "Anders Boström"
foreach ($c in $emptyTextCols) {
    $ws.Cells.Item(4, $c).Formula = '=""'
}

# Row 5
$ws.Cells.Item(5, 1).Value = 112551642
$ws.Cells.Item(5, 2).Value = 77402
$ws.Cells.Item(5, 3).Value = "Ovaliderad"
$ws.Cells.Item(5, 4).Value = "NT"
$ws.Cells.Item(5, 5).Value = 6446
$ws.Cells.Item(5, 6).Value = "Kolflarnlav"
$ws.Cells.Item(5, 7).Value = "Carbonicola anthracophila"
$ws.Cells.Item(5, 8).Value = "(Nyl.) Bendiksby & Timdal"
$ws.Cells.Item(5, 16).Value = "Sydväst Åskaksdammen, Vrm"
$ws.Cells.Item(5, 17).Value = 409388
$ws.Cells.Item(5, 18).Value = 6717259
$ws.Cells.Item(5, 19).Value = 10
$ws.Cells.Item(5, 20).Value = "Värmland"
$ws.Cells.Item(5, 21).Value = "Torsby"
$ws.Cells.Item(5, 22).Value = "Värmland"
$ws.Cells.Item(5, 23).Value = "Norra Ny"
$ws.Cells.Item(5, 25).Value = "2022-08-18"
$ws.Cells.Item(5, 27).Value = "2022-08-18"
$ws.Cells.Item(5, 30).Value = $false
$ws.Cells.Item(5, 31).Value = $false
$ws.Cells.Item(5, 33).Value = $false
$ws.Cells.Item(5, 49).Value = "Anders Boström"
$ws.Cells.Item(5, 50).Value = "Anders Boström"
foreach ($c in $emptyTextCols) {
    $ws.Cells.Item(5, $c).Formula = '=""'
}

# Restore the default "Normal" style on the date columns so no stray
# number-format styling is left behind on cells that should render with the
# workbook's default style.
foreach ($r in 3..5) {
    foreach ($c in $dateCols) {
        $ws.Cells.Item($r, $c).Style = "Normal"
    }
}
